# Add the I0 (column I) and IF (column J) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same bold / centered / bordered style
# (style index 1) that the rest of row 1 already uses.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data ----------------------------------------------------------------
# I2:J23, row by row: column I = I0, column J = IF
$values = @(
    @(1, 1),
    @(1, 5),
    @(1, 8),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(8, 8),
    @(6, 7),
    @(8, 9),
    @(5, 6),
    @(7, 8),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 2),
    @(1, 3),
    @(6, 7),
    @(1, 2)
)

for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = 2 + $idx
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$wb.Save()
